$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DoctorDetails")

$ws.Range("A2").Value = "Dr. Balaji K"
$ws.Range("B2").Value = "Dentist"
$ws.Range("C2").Value = "28 years experience overall"
$ws.Range("D2").Value = "Anna Nagar,Chennai  KB DENTAL CLINIC ORTHODONTIC AND IMPLANT CENTER"
$ws.Range("E2").Value = "₹600 Consultation fee at clinic"

$ws.Range("A6").Value = "Dr. Abhilash Bhaskaran"
$ws.Range("B6").Value = "Dentist"
$ws.Range("C6").Value = "27 years experience overall"
$ws.Range("D6").Value = "Perumbakkam,Chennai  Gleneagles Health City + 1 more"
$ws.Range("E6").Value = "₹800 Consultation fee at clinic"
